$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1283040.337999999
$ws.Range("B3").Value = 7771.034000000002
$ws.Range("B4").Value = 928523.1649999979
$ws.Range("B5").Value = 706.423
$ws.Range("B7").Value = 925963.0650000017
